$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# New values for columns D..AQ, shared by rows 2 and 3
$values = @{
    "D"  = 0.233
    "E"  = 0.252
    "F"  = -0.0839
    "G"  = 0.718479769403224
    "H"  = 0.718479769403224
    "I"  = 0.5636350899270569
    "J"  = 0.4745713934867212
    "K"  = 436.4
    "L"  = 0.4658908935625066
    "M"  = 72.72
    "N"  = 0.01497312990302057
    "O"  = 0.1666361136571952
    "P"  = 66.2
    "Q"  = 0.01363065455968044
    "R"  = 0.1516956920256645
    "S"  = 6.519999999999996
    "T"  = 0.0896589658965896
    "U"  = 118.7
    "V"  = 0.02444046368933638
    "W"  = 0.152843933875035
    "X"  = 0.01921858328392602
    "Y"  = 0.133625350591109
    "Z"  = 0.3078955918229607
    "AA" = 0.1461184400598412
    "AB" = 0.02038118342966721
    "AC" = 0.125737256630174
    "AD" = 424.7
    "AE" = 12.56505632662863
    "AF" = 437.2650563266286
    "AG" = 318.5650563266286
    "AH" = 0.08259689130438606
    "AI" = 0.1045147254780482
    "AJ" = 0.06155531221288676
    "AK" = 0.07836653335494168
    "AL" = 9.140000000000001
    "AM" = 9.140000000000001
    "AN" = 0.7953630353765192
    "AO" = 57.71334792122538
    "AP" = 0.5965972925943941
    "AQ" = 57.71334792122538
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
